$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newText = @'
questions = [
    {
        "title": "You are tasked with creating an Internet of Things (IoT) temperature sensor that operates wirelessly and has low power consumption. This sensor should be able to run on a coin cell battery and activate an alarm within 10 seconds of detecting a temperature that exceeds a specific threshold.Which of the following microcontroller options should you use?",
        "ques_type": 2,
        "options": [
            "A microcontroller that supports a real-time operating system (RTOS).",
            "A microcontroller that supports multiple performance modes with Bluetooth and Wi-Fi connectivity.",
            "A microcontroller that supports multiple power modes with Bluetooth connectivity.",
            "A microcontroller that supports multiple power modes with Wi-Fi connectivity."
        ],
        "score": "A microcontroller that supports multiple power modes with Bluetooth connectivity."
    },
    {
        "title": "You work for a company that develops game consoles. Your company is developing a wireless four-button game controller that is designed to go to sleep to save battery life when there is no activity for one minute. It wakes up from sleep when any button is pressed. The schematic below is provided to you by the hardware engineer.How would you configure the GPIOs in the MCU before it goes to sleep to implement the feature?",
        "ques_type": 2,
        "options": [
            "GPIOs as input, pull up, sensing rising edge.",
            "GPIOs as input, pull up, sensing falling edge.",
            "GPIOs as input, pull down, sensing rising edge.",
            "GPIOs as input, pull down, sensing falling edge."
        ],
        "score": "GPIOs as input, pull down, sensing rising edge."
    },
    {
        "title": "You are asked to debug the code for an ambient light sensor (shown below) which was written by your colleague. The UART on MCU is connected to the PC via a USB-to-UART bridge. You are not receiving messages in the UART console on your PC. Your colleague suggests two possible issues: Either the sensor value (AmbientLightIntensity) is incorrect, or the UART console is not configured correctly.Which step should you take first to find the problem?",
        "ques_type": 2,
        "options": [
            "Measuring the light sensor voltage using a multimeter.",
            "Analyzing UART signals on an oscilloscope.",
            "Setting a breakpoint using a debugger.",
            "Configuring a different UART baud rate and checking messages on the console."
        ],
        "score": "Setting a breakpoint using a debugger."
    },
    {
        "title": "You are asked to write a program for a heart rate monitor (HRM) that uses an ECG (electrocardiography) sensor. The ECG sensor needs to be read every 20 milliseconds to estimate the heart rate correctly. You wrote the program shown below, but the displayed heart rate is wrong. On debugging, you found that the ReadEcgSensor function is executed around every 30 milliseconds.What should you change in the code to display the correct heart rate?",
        "ques_type": 2,
        "options": [
            "Set up a 20-millisecond timer interrupt. Call ReadEcgSensor and EstimateHeartRate functions inside the interrupt service routine.",
            "Set up a 20-millisecond timer interrupt. Call ReadEcgSensor and UpdateDisplay functions inside the interrupt service routine.",
            "Set up a 20-millisecond timer interrupt. Call only ReadEcgSensor inside the interrupt service routine. Update the display more frequently using blocking APIs.",
            "Set up a 20-millisecond timer interrupt. Call only ReadEcgSensor inside the interrupt service routine. Update the display less frequently using non-blocking APIs. "
        ],
        "score": "Set up a 20-millisecond timer interrupt. Call only ReadEcgSensor inside the interrupt service routine. Update the display less frequently using non-blocking APIs."
    }
]
'@

# Remove the old content in A2 entirely (it will no longer exist; the sheet
# dimension collapses back down to just A1).
$ws.Range("A2").ClearContents()

# A1 should hold the full questions text (as a shared string) without any
# special bold/bordered/centered styling - reset it back to the default style.
$ws.Range("A1").ClearFormats()
$ws.Range("A1").Value = $newText
